$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the 4 query scripts (recurrence score filter 31-35 -> 16-20) ---
# CasesTab row (row 2): query (B2) + StatQuery (C2)
$ws.Cells.Item(2, 2).Value = 'MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"] and d.tumor_size_group In ["(2,3]"] and tp.chemotherapy_regimen In ["IV CMF (3 week cycles)"]
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`'
$ws.Cells.Item(2, 3).Value = 'MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"] and d.tumor_size_group In ["(2,3]"] and tp.chemotherapy_regimen In ["IV CMF (3 week cycles)"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files'

# SamplesTab row (row 3): query (B3) + StatQuery (C3)
$ws.Cells.Item(3, 2).Value = 'MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"] and d.tumor_size_group In ["(2,3]"] and tp.chemotherapy_regimen In ["IV CMF (3 week cycles)"]
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`'
$ws.Cells.Item(3, 3).Value = 'MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"]and d.tumor_size_group In ["(2,3]"] and tp.chemotherapy_regimen In ["IV CMF (3 week cycles)"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files'

# FilesTab row (row 4): query (B4) + StatQuery (C4)
$ws.Cells.Item(4, 2).Value = 'MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"] and d.tumor_size_group In ["(2,3]"] and tp.chemotherapy_regimen In ["IV CMF (3 week cycles)"]
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`'
$ws.Cells.Item(4, 3).Value = 'MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"] and d.tumor_size_group In ["(2,3]"] and tp.chemotherapy_regimen In ["IV CMF (3 week cycles)"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files'

# --- Row heights for the (now taller) wrapped query rows ---
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 360
$ws.Rows.Item(4).RowHeight = 360

# --- Column width nudges (best effort; engine quantizes to 1/6 char) ---
$ws.Columns.Item(1).ColumnWidth = 10.022135416666666
$ws.Columns.Item(2).ColumnWidth = 75.02213541666667
$ws.Columns.Item(3).ColumnWidth = 59.592447916666664
$ws.Columns.Item(4).ColumnWidth = 69.45182291666667

# --- Scroll / selection: view now parked on row 4 ---
$ws.Range("C4").Select()
